$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 6456.25
$ws.Range("I74").Value = 7920
$ws.Range("J74").Value = 4016.6667
$ws.Range("K74").Value = 7920
$ws.Range("L74").Value = 4016.6667
$ws.Range("M74").Value = -6984
$ws.Range("N74").Value = -5888.6667

$ws.Range("H76").Value = 5519.7
$ws.Range("I76").Value = 4986.625
$ws.Range("K76").Value = 4986.625
$ws.Range("M76").Value = -4671.625

$ws.Range("H77").Value = 6456.25
$ws.Range("I77").Value = 7920
$ws.Range("J77").Value = 4016.6667
$ws.Range("K77").Value = 39600
$ws.Range("L77").Value = 20083.3335
$ws.Range("M77").Value = -34920
$ws.Range("N77").Value = -29443.3335

$ws.Range("H79").Value = 5519.7
$ws.Range("I79").Value = 4986.625
$ws.Range("K79").Value = 4986.625
$ws.Range("M79").Value = -3894.625

$ws.Range("H80").Value = 622.0833
$ws.Range("I80").Value = 496.81818
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 1490.45454
$ws.Range("L80").Value = 6000
$ws.Range("M80").Value = -492.45454
$ws.Range("N80").Value = -7996

$ws.Range("H83").Value = 622.0833
$ws.Range("I83").Value = 496.81818
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 4471.36362
$ws.Range("L83").Value = 18000
$ws.Range("M83").Value = 520.6363799999999
$ws.Range("N83").Value = -27984

$ws.Range("H138").Value = 2737.4949
$ws.Range("I138").Value = 1489.25
$ws.Range("J138").Value = 3352.0154
$ws.Range("K138").Value = 4467.75
$ws.Range("L138").Value = 10056.0462
$ws.Range("M138").Value = 672.25
$ws.Range("N138").Value = -20336.0462

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17219.393
$ws.Range("I32").Value = 15374.707
$ws.Range("K32").Value = 15374.707
$ws.Range("M32").Value = -15087.707

$ws.Range("H61").Value = 5177.763
$ws.Range("I61").Value = 5266.484
$ws.Range("K61").Value = 5266.484
$ws.Range("M61").Value = -5054.484

$ws.Range("H63").Value = 4169803.5
$ws.Range("I63").Value = 11112542
$ws.Range("J63").Value = 4160
$ws.Range("K63").Value = 11112542
$ws.Range("L63").Value = 4160
$ws.Range("M63").Value = -11111856
$ws.Range("N63").Value = -5532

$ws.Range("H66").Value = 4169803.5
$ws.Range("I66").Value = 11112542
$ws.Range("J66").Value = 4160
$ws.Range("K66").Value = 55562710
$ws.Range("L66").Value = 20800
$ws.Range("M66").Value = -55559278
$ws.Range("N66").Value = -27664

$ws.Range("H74").Value = 658.5472
$ws.Range("I74").Value = 517.6857
$ws.Range("J74").Value = 932.44446
$ws.Range("K74").Value = 517.6857
$ws.Range("L74").Value = 932.44446
$ws.Range("M74").Value = 356.3143
$ws.Range("N74").Value = -2680.44446

$ws.Range("H77").Value = 658.5472
$ws.Range("I77").Value = 517.6857
$ws.Range("J77").Value = 932.44446
$ws.Range("K77").Value = 2588.4285
$ws.Range("L77").Value = 4662.2223
$ws.Range("M77").Value = 1779.5715
$ws.Range("N77").Value = -13398.2223

$ws.Range("H136").Value = 5177.763
$ws.Range("I136").Value = 5266.484
$ws.Range("K136").Value = 15799.452
$ws.Range("M136").Value = -13249.452

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3045
$ws.Range("I86").Value = 2351.25
$ws.Range("K86").Value = 2351.25
$ws.Range("M86").Value = -1228.25

$ws.Range("H89").Value = 3045
$ws.Range("I89").Value = 2351.25
$ws.Range("K89").Value = 11756.25
$ws.Range("M89").Value = -6140.25

$ws.Range("H105").Value = 2977.7778
$ws.Range("I105").Value = 4000
$ws.Range("J105").Value = 2160
$ws.Range("K105").Value = 4000
$ws.Range("L105").Value = 2160
$ws.Range("M105").Value = -2253
$ws.Range("N105").Value = -5654

$ws.Range("H134").Value = 7258379.5
$ws.Range("I134").Value = 7949313.5
$ws.Range("J134").Value = 3571
$ws.Range("K134").Value = 23847940.5
$ws.Range("L134").Value = 10713
$ws.Range("M134").Value = -23845405.5
$ws.Range("N134").Value = -15783

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10013.275
$ws.Range("I31").Value = 1230.9375
$ws.Range("J31").Value = 20822.309
$ws.Range("K31").Value = 1230.9375
$ws.Range("L31").Value = 20822.309
$ws.Range("M31").Value = -935.9375
$ws.Range("N31").Value = -21412.309

$ws.Range("H34").Value = 10013.275
$ws.Range("I34").Value = 1230.9375
$ws.Range("J34").Value = 20822.309
$ws.Range("K34").Value = 1230.9375
$ws.Range("L34").Value = 20822.309
$ws.Range("M34").Value = -1028.9375
$ws.Range("N34").Value = -21226.309

$ws.Range("H58").Value = 3599819.8
$ws.Range("I58").Value = 4795989.5
$ws.Range("J58").Value = 11310.4
$ws.Range("K58").Value = 4795989.5
$ws.Range("L58").Value = 11310.4
$ws.Range("M58").Value = -4795786.5
$ws.Range("N58").Value = -11716.4

$ws.Range("H107").Value = 303.96774
$ws.Range("I107").Value = 330.8846
$ws.Range("J107").Value = 164
$ws.Range("K107").Value = 330.8846
$ws.Range("L107").Value = 164
$ws.Range("M107").Value = 1589.1154
$ws.Range("N107").Value = -4004

$ws.Range("H132").Value = 5293902.5
$ws.Range("I132").Value = 7752795
$ws.Range("J132").Value = 7283.5
$ws.Range("K132").Value = 23258385
$ws.Range("L132").Value = 21850.5
$ws.Range("M132").Value = -23255855
$ws.Range("N132").Value = -26910.5

$ws.Range("H134").Value = 10969249
$ws.Range("I134").Value = 13890803
$ws.Range("J134").Value = 7212965
$ws.Range("K134").Value = 41672409
$ws.Range("L134").Value = 21638895
$ws.Range("M134").Value = -41669874
$ws.Range("N134").Value = -21643965

$ws.Range("H136").Value = 3599819.8
$ws.Range("I136").Value = 4795989.5
$ws.Range("J136").Value = 11310.4
$ws.Range("K136").Value = 14387968.5
$ws.Range("L136").Value = 33931.2
$ws.Range("M136").Value = -14385418.5
$ws.Range("N136").Value = -39031.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 7150446.5
$ws.Range("I113").Value = 467.25
$ws.Range("J113").Value = 10010438
$ws.Range("K113").Value = 1401.75
$ws.Range("L113").Value = 30031314
$ws.Range("M113").Value = 768.25
$ws.Range("N113").Value = -30035654

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 26120.133
$ws.Range("I70").Value = 35049.688
$ws.Range("J70").Value = 4139.6924
$ws.Range("K70").Value = 35049.688
$ws.Range("L70").Value = 4139.6924
$ws.Range("M70").Value = -34779.688
$ws.Range("N70").Value = -4679.6924

$ws.Range("H73").Value = 26120.133
$ws.Range("I73").Value = 35049.688
$ws.Range("J73").Value = 4139.6924
$ws.Range("K73").Value = 35049.688
$ws.Range("L73").Value = 4139.6924
$ws.Range("M73").Value = -34113.688
$ws.Range("N73").Value = -6011.6924

$ws.Range("H80").Value = 2816.6667
$ws.Range("J80").Value = 3133.3333
$ws.Range("L80").Value = 3133.3333
$ws.Range("N80").Value = -5129.3333

$ws.Range("H83").Value = 2816.6667
$ws.Range("J83").Value = 3133.3333
$ws.Range("L83").Value = 15666.6665
$ws.Range("N83").Value = -25650.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4279043
$ws.Range("I132").Value = 4904564.5
$ws.Range("K132").Value = 14713693.5
$ws.Range("M132").Value = -14711163.5

$ws.Range("H136").Value = 3590.4807
$ws.Range("I136").Value = 4838.125
$ws.Range("J136").Value = 1594.25
$ws.Range("K136").Value = 14514.375
$ws.Range("L136").Value = 4782.75
$ws.Range("M136").Value = -11964.375
$ws.Range("N136").Value = -9882.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 7143588.5
$ws.Range("I126").Value = 8621349
$ws.Range("J126").Value = 1079.6666
$ws.Range("K126").Value = 25864047
$ws.Range("L126").Value = 3238.9998
$ws.Range("M126").Value = -25861577
$ws.Range("N126").Value = -8178.9998

$ws.Range("H132").Value = 194972.36
$ws.Range("I132").Value = 14643.62
$ws.Range("J132").Value = 2328862.5
$ws.Range("K132").Value = 43930.86
$ws.Range("L132").Value = 6986587.5
$ws.Range("M132").Value = -41400.86
$ws.Range("N132").Value = -6991647.5

$ws.Range("H136").Value = 417498.84
$ws.Range("I136").Value = 495481.66
$ws.Range("J136").Value = 1590.5555
$ws.Range("K136").Value = 1486444.98
$ws.Range("L136").Value = 4771.666499999999
$ws.Range("M136").Value = -1483894.98
$ws.Range("N136").Value = -9871.666499999999
